$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 1504
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("H43").Value = 1491.25
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 1535.909
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 1535.909
$ws.Range("N43").Value = -1673.909
$ws.Range("M43").Value = -931
$ws.Range("H106").Value = 3599.8
$ws.Range("I106").Value = 2544.0908
$ws.Range("K106").Value = 2544.0908
$ws.Range("M106").Value = -1913.0908
$ws.Range("H137").Value = 1471
$ws.Range("I137").Value = 894.5
$ws.Range("J137").Value = 2212.2144
$ws.Range("K137").Value = 2683.5
$ws.Range("L137").Value = 6636.6432
$ws.Range("M137").Value = -133.5
$ws.Range("N137").Value = -11736.6432
$ws.Range("H138").Value = 2492.054
$ws.Range("J138").Value = 2054.0715
$ws.Range("L138").Value = 6162.2145
$ws.Range("N138").Value = -16442.2145
$ws.Range("N10").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3864.2373
$ws.Range("I32").Value = 2480.3257
$ws.Range("K32").Value = 2480.3257
$ws.Range("M32").Value = -2193.3257
$ws.Range("H45").Value = 1668.375
$ws.Range("I45").Value = 996.5
$ws.Range("J45").Value = 1892.3334
$ws.Range("K45").Value = 996.5
$ws.Range("L45").Value = 1892.3334
$ws.Range("N45").Value = -2646.3334
$ws.Range("M45").Value = -619.5
$ws.Range("H61").Value = 3484.0476
$ws.Range("I61").Value = 2252.5625
$ws.Range("K61").Value = 2252.5625
$ws.Range("M61").Value = -2040.5625
$ws.Range("H122").Value = 1630.2941
$ws.Range("I122").Value = 1366.1111
$ws.Range("K122").Value = 4098.3333
$ws.Range("M122").Value = -1648.3333
$ws.Range("H132").Value = 1657.7858
$ws.Range("I132").Value = 1100.9166
$ws.Range("K132").Value = 3302.7498
$ws.Range("M132").Value = -772.7498000000001
$ws.Range("H136").Value = 3484.0476
$ws.Range("I136").Value = 2252.5625
$ws.Range("K136").Value = 6757.6875
$ws.Range("M136").Value = -4207.6875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 8047.95
$ws.Range("I134").Value = 8817.5
$ws.Range("J134").Value = 5739.3
$ws.Range("K134").Value = 26452.5
$ws.Range("L134").Value = 17217.9
$ws.Range("M134").Value = -23917.5
$ws.Range("N134").Value = -22287.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2805.55
$ws.Range("I31").Value = 1379.3572
$ws.Range("K31").Value = 1379.3572
$ws.Range("M31").Value = -1084.3572
$ws.Range("H34").Value = 2805.55
$ws.Range("I34").Value = 1379.3572
$ws.Range("K34").Value = 1379.3572
$ws.Range("M34").Value = -1177.3572
$ws.Range("H53").Value = 60000
$ws.Range("J53").Value = 60000
$ws.Range("L53").Value = 60000
$ws.Range("N53").Value = -61214

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 3000
$ws.Range("J116").Value = 3000
$ws.Range("L116").Value = 9000
$ws.Range("N116").Value = -15884
$ws.Range("H131").Value = 777.78
$ws.Range("I131").Value = 355.125
$ws.Range("J131").Value = 814.5326
$ws.Range("K131").Value = 1065.375
$ws.Range("L131").Value = 2443.5978
$ws.Range("M131").Value = 3974.625
$ws.Range("N131").Value = -12523.5978
$ws.Range("H140").Value = 1702.4517
$ws.Range("I140").Value = 855.7143
$ws.Range("K140").Value = 2567.1429
$ws.Range("M140").Value = 2612.8571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2445.5833
$ws.Range("I122").Value = 2271.5
$ws.Range("J122").Value = 2793.75
$ws.Range("K122").Value = 6814.5
$ws.Range("L122").Value = 8381.25
$ws.Range("M122").Value = -4364.5
$ws.Range("N122").Value = -13281.25
$ws.Range("H132").Value = 2964298.5
$ws.Range("I132").Value = 4277726.5
$ws.Range("J132").Value = 9085.75
$ws.Range("K132").Value = 12833179.5
$ws.Range("L132").Value = 27257.25
$ws.Range("M132").Value = -12830649.5
$ws.Range("N132").Value = -32317.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3829.4546
$ws.Range("I7").Value = 1738.5714
$ws.Range("J7").Value = 7488.5
$ws.Range("K7").Value = 1738.5714
$ws.Range("L7").Value = 7488.5
$ws.Range("M7").Value = -1626.5714
$ws.Range("N7").Value = -7712.5
$ws.Range("H16").Value = 2370.1765
$ws.Range("I16").Value = 3983.889
$ws.Range("J16").Value = 554.75
$ws.Range("K16").Value = 3983.889
$ws.Range("L16").Value = 554.75
$ws.Range("M16").Value = -3813.889
$ws.Range("N16").Value = -894.75
$ws.Range("H55").Value = 525.125
$ws.Range("J55").Value = 548.4
$ws.Range("L55").Value = 548.4
$ws.Range("N55").Value = -894.4
$ws.Range("H64").Value = 512399.5
$ws.Range("J64").Value = 24800
$ws.Range("L64").Value = 24800
$ws.Range("N64").Value = -25250
$ws.Range("H67").Value = 512399.5
$ws.Range("J67").Value = 24800
$ws.Range("L67").Value = 24800
$ws.Range("N67").Value = -26360
$ws.Range("H94").Value = 47693.332
$ws.Range("J94").Value = 47693.332
$ws.Range("L94").Value = 47693.332
$ws.Range("N94").Value = -49045.332
$ws.Range("H100").Value = 1800
$ws.Range("I100").Value = 1800
$ws.Range("K100").Value = 1800
$ws.Range("M100").Value = -1259
$ws.Range("H122").Value = 6618.1816
$ws.Range("I122").Value = 5375
$ws.Range("J122").Value = 7328.5713
$ws.Range("K122").Value = 16125
$ws.Range("L122").Value = 21985.7139
$ws.Range("M122").Value = -13675
$ws.Range("N122").Value = -26885.7139
$ws.Range("H126").Value = 3829.4546
$ws.Range("I126").Value = 1738.5714
$ws.Range("J126").Value = 7488.5
$ws.Range("K126").Value = 5215.7142
$ws.Range("L126").Value = 22465.5
$ws.Range("M126").Value = -2745.7142
$ws.Range("N126").Value = -27405.5
$ws.Range("H132").Value = 2268.8
$ws.Range("I132").Value = 1336.25
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 4008.75
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -1478.75
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 1259306.2
$ws.Range("J63").Value = 1676999.6
$ws.Range("L63").Value = 1676999.6
$ws.Range("N63").Value = -1678247.6
$ws.Range("H66").Value = 1259306.2
$ws.Range("J66").Value = 1676999.6
$ws.Range("L66").Value = 5030998.800000001
$ws.Range("N66").Value = -5037238.800000001
$ws.Range("H95").Value = 49999.5
$ws.Range("J95").Value = 49999.5
$ws.Range("L95").Value = 49999.5
$ws.Range("N95").Value = -55491.5
$ws.Range("H126").Value = 4504.375
$ws.Range("I126").Value = 3924.5938
$ws.Range("K126").Value = 11773.7814
$ws.Range("M126").Value = -9303.7814
$ws.Range("H132").Value = 8256.040000000001
$ws.Range("I132").Value = 2801.3333
$ws.Range("K132").Value = 8403.999899999999
$ws.Range("M132").Value = -8403.999899999999

